$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Marks")

# Session 7 (Branch&Bound) mark for the student
$ws.Range("H4").Value = 9.5

# Comment for Session 7 (Branch&Bound), matching the feedback row (row 5)
# (set before I4 so the shared-string table gets this entry first, at index 24)
$ws.Range("H5").Value = "Very good. To greatly improve times, you may also comment lines 35 and 37 in Heap.java file. That part is only used to check if a node is already repeated but in this problem that will never happen. The times are going to be much much better without it. Be careful with the toString() method. In order to recoginize two nodes as different, the toString() method should print different information. In your case, two nodes are consider equals if the have the same songs in the same blocks, but if they are at different levels of the tree they should be considered as different. To sum up: you also need to print the level of the node in the toString() method."

# Test mark column -> "NA"
$ws.Range("I4").Value = "NA"

# Scroll/selection state as left by the author (scrolled one column right,
# so column B is the left-most visible column, with I5:I12 selected)
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("I5:I12").Select()
